$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "index" row (currently row 2). This shifts every row below it up by one.
$ws.Rows.Item(2).Delete()

# Re-write the resulting data block (rows 2-12) with the refreshed category
# ranking values used by the new "gabarito" generator.
$data = @(
    @(10, "Ignore", 690),
    @(9, "Other", 46),
    @(7, "time_manipulation", 23),
    @(0, "access_control", 19),
    @(3, "reentrancy", 16),
    @(4, "unchecked_low_calls", 11),
    @(2, "denial_service", 3),
    @(1, "arithmetic", 0),
    @(5, "bad_randomness", 0),
    @(6, "front_running", 0),
    @(8, "short_addresses", 0)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
